$d = $word.ActiveDocument

# --- 1) Insert ", nuestros nombres son: Nicolás, Daniela y Sebastián." right
#        after "Hola" (end of the existing run), as its own run carrying the
#        same es-ES language formatting. ---
$helloEnd = $d.Paragraphs.Item(1).Range.Start + 4   # "Hola" is 4 chars
$insert1 = $d.Range($helloEnd, $helloEnd)
$namesText = ", nuestros nombres son: Nicolás, Daniela y Sebastián."
$insert1.InsertAfter($namesText)

$namesRange = $d.Range($helloEnd, $helloEnd + $namesText.Length)
$namesRange.LanguageID = "es-ES"
# Nudge formatting so the engine materializes this as its own <w:r>
# instead of silently re-merging it into the preceding "Hola" run.
$namesRange.Bold = 1
$namesRange.Bold = 0

# --- 2) Insert a trailing space as a new run placed *after* the
#        _GoBack bookmark (bookmarkStart/bookmarkEnd), still inside the
#        same paragraph, before the paragraph mark. ---
$paraEnd = $d.Paragraphs.Item(1).Range.End - 1   # position right before the pilcrow, after the bookmark
$insert2 = $d.Range($paraEnd, $paraEnd)
$insert2.InsertAfter(" ")

$spaceRange = $d.Range($paraEnd, $paraEnd + 1)
$spaceRange.LanguageID = "es-ES"
$spaceRange.Bold = 1
$spaceRange.Bold = 0
